$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 already has A18 = "017" (id). Add the rest of the simulation row data,
# matching the pattern used by the other rows (e.g. row 11, which also has E=1000).
$ws.Range("B18").Value = "1 / 1.8"
$ws.Range("C18").Value = 0.3
$ws.Range("D18").Value = 1.8
$ws.Range("E18").Value = 1000
$ws.Range("F18").Value = "children, adolescents, adults, elderly"
$ws.Range("G18").Value = "open, close"
$ws.Range("H18").Value = 43
$ws.Range("I18").Value = 0.5
$ws.Range("J18").Value = $false
$ws.Range("K18").Value = 1000000

# Copy the style of the corresponding cells in row 11 (same data pattern) down to row 18
$ws.Range("B11:D11").Copy()
$ws.Range("B18:D18").PasteSpecial(-4122)
$ws.Range("E11").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("F11:K11").Copy()
$ws.Range("F18:K18").PasteSpecial(-4122)

# Update the view: frozen-pane top-left cell and active selection
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("N28").Select()
